# API: Gameweeks import (#25)
# Adds two new columns to the "Challenges" sheet:
#   S = "Show Statistics Continuously" (header) / "true" (row 2, as literal text)
#   T = "Gameweek" (header) / 1 (row 2, numeric)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# New header cells
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# New numeric value - plain number, no special handling needed
$ws.Range("T2").Value = 1

# New "true" value must be stored as literal text (shared string), not as an
# Excel boolean. Typing the literal word true/false straight into .Value
# triggers Excel's own boolean auto-detection, so instead we build it via a
# formula that evaluates to the text "true" and then convert that formula to
# a static value in place (copy / paste-values onto itself) - this keeps the
# cell a plain inline text value with no residual formula and no formatting
# changes.
$ws.Range("S2").Formula = '=""&"true"'
$ws.Range("S2").Copy()
$ws.Range("S2").PasteSpecial(-4163)
$excel.CutCopyMode = $false
